$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "27.454.30"

Set-TextValue $ws.Range("D3") "1.822.81"
Set-TextValue $ws.Range("E3") "  -0.99%  "

Set-TextValue $ws.Range("E4") "  +0.01%  "

Set-TextValue $ws.Range("D5") "312.30"
Set-TextValue $ws.Range("E5") "  +0.14%  "

Set-TextValue $ws.Range("E6") "  -0.01%  "

Set-TextValue $ws.Range("D7") "0.4238"
Set-TextValue $ws.Range("E7") "  -0.48%  "

Set-TextValue $ws.Range("D8") "0.3622"
Set-TextValue $ws.Range("E8") "  +0.48%  "

Set-TextValue $ws.Range("D9") "0.07193"
Set-TextValue $ws.Range("E9") "  -1.28%  "

Set-TextValue $ws.Range("E10") "  -0.88%  "

Set-TextValue $ws.Range("D11") "20.59"
Set-TextValue $ws.Range("E11") "  +0.03%  "

Set-TextValue $ws.Range("D12") "1.878.96"
Set-TextValue $ws.Range("E12") "  +2.73%  "

Set-TextValue $ws.Range("D13") "5.393"

Set-TextValue $ws.Range("D14") "6.465"
Set-TextValue $ws.Range("E14") "  -0.45%  "

Set-TextValue $ws.Range("D15") "0.06925"
Set-TextValue $ws.Range("E15") "  -0.48%  "

Set-TextValue $ws.Range("D17") "80.26"
Set-TextValue $ws.Range("E17") "  +1.22%  "

Set-TextValue $ws.Range("D18") "0.000008890"
Set-TextValue $ws.Range("E18") "  -0.37%  "

Set-TextValue $ws.Range("E19") "  +0.02%  "

Set-TextValue $ws.Range("D20") "15.33"
Set-TextValue $ws.Range("E20") "  +0.73%  "

Set-TextValue $ws.Range("D21") "27.384.01"
Set-TextValue $ws.Range("E21") "  -0.85%  "

Set-TextValue $ws.Range("D22") "5.135"
Set-TextValue $ws.Range("E22") "  +3.42%  "

Set-TextValue $ws.Range("D23") "10.92"
Set-TextValue $ws.Range("E23") "  +5.86%  "

Set-TextValue $ws.Range("D24") "2.034.86"
Set-TextValue $ws.Range("E24") "  -1.64%  "

Set-TextValue $ws.Range("D25") "1.983"
Set-TextValue $ws.Range("E25") "  +0.28%  "

Set-TextValue $ws.Range("D26") "154.89"
Set-TextValue $ws.Range("E26") "  -0.31%  "

Set-TextValue $ws.Range("D27") "18.68"
Set-TextValue $ws.Range("E27") "  +1.17%  "

Set-TextValue $ws.Range("D28") "5.154"
Set-TextValue $ws.Range("E28") "  -0.87%  "

Set-TextValue $ws.Range("D29") "113.97"
Set-TextValue $ws.Range("E29") "  -4.66%  "

Set-TextValue $ws.Range("D30") "1.791"
Set-TextValue $ws.Range("E30") "  -3.78%  "

Set-TextValue $ws.Range("D31") "0.08836"
Set-TextValue $ws.Range("E31") "  -0.50%  "

Set-TextValue $ws.Range("D32") "0.7498"
Set-TextValue $ws.Range("E32") "  -1.32%  "

Set-TextValue $ws.Range("D33") "2.969"
Set-TextValue $ws.Range("E33") "  +0.13%  "

Set-TextValue $ws.Range("D34") "4.531"
Set-TextValue $ws.Range("E34") "  +1.19%  "

Set-TextValue $ws.Range("D35") "1.118"
Set-TextValue $ws.Range("E35") "  -0.39%  "

Set-TextValue $ws.Range("E36") "  +0.01%  "

Set-TextValue $ws.Range("D37") "1.088"
Set-TextValue $ws.Range("E37") "  -1.54%  "

Set-TextValue $ws.Range("D38") "0.05281"
Set-TextValue $ws.Range("E38") "  -2.39%  "

Set-TextValue $ws.Range("D39") "0.01916"
Set-TextValue $ws.Range("E39") "  -0.33%  "

Set-TextValue $ws.Range("D40") "2.777"
Set-TextValue $ws.Range("E40") "  -1.25%  "

Set-TextValue $ws.Range("D41") "0.5059"
Set-TextValue $ws.Range("E41") "  +0.27%  "

Set-TextValue $ws.Range("D42") "0.1641"
Set-TextValue $ws.Range("E42") "  -0.82%  "

Set-TextValue $ws.Range("D43") "6.447"
Set-TextValue $ws.Range("E43") "  -1.18%  "

Set-TextValue $ws.Range("D44") "8.327"
Set-TextValue $ws.Range("E44") "  -0.23%  "

Set-TextValue $ws.Range("D45") "10.40"
Set-TextValue $ws.Range("E45") "  +0.47%  "

Set-TextValue $ws.Range("D46") "105.59"
Set-TextValue $ws.Range("E46") "  -0.28%  "

Set-TextValue $ws.Range("B47") "Decentraland"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D47") "0.4675"
Set-TextValue $ws.Range("E47") "  +1.31%  "

Set-TextValue $ws.Range("B48") "Cronos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.06444"
Set-TextValue $ws.Range("E48") "  -1.41%  "

Set-TextValue $ws.Range("E49") "  -0.05%  "

Set-TextValue $ws.Range("E50") "  -1.04%  "

Set-TextValue $ws.Range("D51") "63.67"
Set-TextValue $ws.Range("E51") "  -0.82%  "
